$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I (I0) and J (IF)
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Match the header style used by the existing header cells (bold, bordered, centered)
# by copying the format from H1, reusing the same style index as the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data for I2:J78 (I0 and IF values per row)
$data = @(
    @{Row=2; I=8; J=8},
    @{Row=3; I=7; J=7},
    @{Row=4; I=8; J=8},
    @{Row=5; I=7; J=7},
    @{Row=6; I=7; J=7},
    @{Row=7; I=8; J=8},
    @{Row=8; I=7; J=7},
    @{Row=9; I=7; J=7},
    @{Row=10; I=7; J=7},
    @{Row=11; I=7; J=7},
    @{Row=12; I=7; J=7},
    @{Row=13; I=8; J=8},
    @{Row=14; I=6; J=7},
    @{Row=15; I=9; J=9},
    @{Row=16; I=8; J=8},
    @{Row=17; I=7; J=7},
    @{Row=18; I=9; J=9},
    @{Row=19; I=7; J=7},
    @{Row=20; I=8; J=8},
    @{Row=21; I=7; J=7},
    @{Row=22; I=7; J=7},
    @{Row=23; I=8; J=8},
    @{Row=24; I=8; J=8},
    @{Row=25; I=8; J=9},
    @{Row=26; I=5; J=6},
    @{Row=27; I=7; J=8},
    @{Row=28; I=7; J=7},
    @{Row=29; I=6; J=6},
    @{Row=30; I=7; J=8},
    @{Row=31; I=10; J=10},
    @{Row=32; I=7; J=7},
    @{Row=33; I=9; J=9},
    @{Row=34; I=8; J=8},
    @{Row=35; I=8; J=9},
    @{Row=36; I=8; J=8},
    @{Row=37; I=8; J=8},
    @{Row=38; I=6; J=6},
    @{Row=39; I=7; J=7},
    @{Row=40; I=7; J=8},
    @{Row=41; I=8; J=8},
    @{Row=42; I=10; J=10},
    @{Row=43; I=7; J=7},
    @{Row=44; I=6; J=7},
    @{Row=45; I=10; J=10},
    @{Row=46; I=8; J=8},
    @{Row=47; I=8; J=8},
    @{Row=48; I=7; J=7},
    @{Row=49; I=7; J=7},
    @{Row=50; I=8; J=9},
    @{Row=51; I=7; J=7},
    @{Row=52; I=9; J=9},
    @{Row=53; I=8; J=8},
    @{Row=54; I=9; J=9},
    @{Row=55; I=8; J=8},
    @{Row=56; I=8; J=8},
    @{Row=57; I=7; J=8},
    @{Row=58; I=8; J=8},
    @{Row=59; I=9; J=9},
    @{Row=60; I=9; J=9},
    @{Row=61; I=8; J=9},
    @{Row=62; I=9; J=9},
    @{Row=63; I=7; J=7},
    @{Row=64; I=9; J=9},
    @{Row=65; I=8; J=9},
    @{Row=66; I=9; J=9},
    @{Row=67; I=7; J=7},
    @{Row=68; I=9; J=9},
    @{Row=69; I=9; J=9},
    @{Row=70; I=7; J=8},
    @{Row=71; I=6; J=6},
    @{Row=72; I=8; J=8},
    @{Row=73; I=6; J=6},
    @{Row=74; I=7; J=7},
    @{Row=75; I=5; J=5},
    @{Row=76; I=5; J=5},
    @{Row=77; I=5; J=5},
    @{Row=78; I=3; J=3}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}
